$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $addr, $val) {
    # Force the cell to remain Text-typed (avoids Excel auto-
    # converting numeric-looking strings to Number), then restore
    # the cell's original style so no stray formatting is left behind.
    $range = $ws.Range($addr)
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $val
    $range.Style = $origStyle
}

function Set-PlainCell($ws, $addr, $val) {
    $ws.Range($addr).Value = $val
}

# Row 2
Set-TextCell $ws "D2" '66.246.06'
Set-PlainCell $ws "E2" '  -1.88%  '

# Row 3
Set-TextCell $ws "D3" '3.432.94'
Set-PlainCell $ws "E3" '  -4.43%  '

# Row 4
Set-PlainCell $ws "E4" '  +0.02%  '

# Row 5
Set-TextCell $ws "D5" '186.94'
Set-PlainCell $ws "E5" '  -6.30%  '

# Row 6
Set-TextCell $ws "D6" '539.82'
Set-PlainCell $ws "E6" '  -3.13%  '

# Row 7
Set-PlainCell $ws "E7" '  +1.10%  '

# Row 8
Set-TextCell $ws "D8" '3.425.72'
Set-PlainCell $ws "E8" '  -4.47%  '

# Row 9
Set-PlainCell $ws "E9" '  -0.10%  '

# Row 10
Set-TextCell $ws "D10" '0.638'
Set-PlainCell $ws "E10" '  -4.59%  '

# Row 11
Set-TextCell $ws "D11" '58.75'
Set-PlainCell $ws "E11" '  -0.60%  '

# Row 12
Set-TextCell $ws "D12" '0.137'
Set-PlainCell $ws "E12" '  -9.36%  '

# Row 13
Set-TextCell $ws "D13" '0.0000260'
Set-PlainCell $ws "E13" '  -9.04%  '

# Row 14
Set-PlainCell $ws "E14" '  -4.44%  '

# Row 15
Set-TextCell $ws "D15" '3.974.15'
Set-PlainCell $ws "E15" '  -4.85%  '

# Row 16
Set-PlainCell $ws "E16" '  -1.96%  '

# Row 17
Set-TextCell $ws "D17" '3.425.10'
Set-PlainCell $ws "E17" '  -4.81%  '

# Row 18
Set-TextCell $ws "D18" '66.014.09'
Set-PlainCell $ws "E18" '  -2.15%  '

# Row 19
Set-TextCell $ws "D19" '17.78'
Set-PlainCell $ws "E19" '  -6.06%  '

# Row 20
Set-TextCell $ws "D20" '11.51'
Set-PlainCell $ws "E20" '  -5.98%  '

# Row 21
Set-TextCell $ws "D21" '0.999'
Set-PlainCell $ws "E21" '  -6.99%  '

# Row 22
Set-TextCell $ws "D22" '386.36'
Set-PlainCell $ws "E22" '  -3.11%  '

# Row 23
Set-TextCell $ws "D23" '83.96'
Set-PlainCell $ws "E23" '  -1.23%  '

# Row 24
Set-PlainCell $ws "E24" '  -7.01%  '

# Row 25
Set-TextCell $ws "D25" '11.15'
Set-PlainCell $ws "E25" '  -13.18%  '

# Row 26
Set-TextCell $ws "D26" '3.79'
Set-PlainCell $ws "E26" '  -1.84%  '

# Row 27
Set-TextCell $ws "D27" '12.01'
Set-PlainCell $ws "E27" '  -3.76%  '

# Row 28
Set-TextCell $ws "D28" '2.74'
Set-PlainCell $ws "E28" '  -6.73%  '

# Row 29
Set-TextCell $ws "D29" '8.71'
Set-PlainCell $ws "E29" '  -7.96%  '

# Row 30
Set-TextCell $ws "D30" '699.57'
Set-PlainCell $ws "E30" '  +5.29%  '

# Row 31
Set-TextCell $ws "D31" '30.22'
Set-PlainCell $ws "E31" '  -3.74%  '

# Row 32
Set-TextCell $ws "D32" '6.91'
Set-PlainCell $ws "E32" '  -17.91%  '

# Row 33
Set-TextCell $ws "D33" '11.42'
Set-PlainCell $ws "E33" '  -6.15%  '

# Row 34
Set-TextCell $ws "D34" '62.27'
Set-PlainCell $ws "E34" '  -2.17%  '

# Row 35
Set-PlainCell $ws "E35" '  -4.30%  '

# Row 36
Set-PlainCell $ws "E36" '  +0.05%  '

# Row 37
Set-TextCell $ws "D37" '37.33'
Set-PlainCell $ws "E37" '  -11.64%  '

# Row 38
Set-PlainCell $ws "E38" '  -9.30%  '

# Row 39
Set-TextCell $ws "D39" '0.999'
Set-PlainCell $ws "E39" '  -0.01%  '

# Row 40
Set-PlainCell $ws "E40" '  -5.86%  '

# Row 41
Set-TextCell $ws "D41" '2.938.82'
Set-PlainCell $ws "E41" '  -9.62%  '

# Row 42
Set-TextCell $ws "D42" '2.84'
Set-PlainCell $ws "E42" '  -11.26%  '

# Row 43
Set-TextCell $ws "D43" '2.74'
Set-PlainCell $ws "E43" '  +0.21%  '

# Row 44
Set-TextCell $ws "D44" ([string]::Concat('0.0', [char]0x2083, '0640'))
Set-PlainCell $ws "E44" '  -16.74%  '

# Row 45
Set-TextCell $ws "D45" '2.43'
Set-PlainCell $ws "E45" '  -13.44%  '

# Row 46
Set-TextCell $ws "D46" '0.0395'
Set-PlainCell $ws "E46" '  -5.03%  '

# Row 47
Set-PlainCell $ws "E47" '  -2.09%  '

# Row 48
Set-PlainCell $ws "B48" 'ApeXProtocol'
Set-PlainCell $ws "C48" 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
Set-TextCell $ws "D48" '2.95'
Set-PlainCell $ws "E48" '  -6.01%  '

# Row 49
Set-PlainCell $ws "B49" 'Monero'
Set-PlainCell $ws "C49" 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextCell $ws "D49" '134.17'
Set-PlainCell $ws "E49" '  -3.78%  '

# Row 50
Set-TextCell $ws "D50" '2.62'
Set-PlainCell $ws "E50" '  -4.06%  '

# Row 51
Set-TextCell $ws "D51" '2.37'
Set-PlainCell $ws "E51" '  -21.45%  '


